$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.952.01"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "3.859.58"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "468.91"
$ws.Range("E5").Value = "  +4.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.14"
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.609"
$ws.Range("E7").Value = "  -2.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.710"
$ws.Range("E9").Value = "  -4.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  +3.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000339"
$ws.Range("E11").Value = "  +3.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.83"
$ws.Range("E12").Value = "  -4.37%  "
$ws.Range("D13").Value = "4.539.33"
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.11"
$ws.Range("E14").Value = "  -2.77%  "
$ws.Range("B15").Value = "Uniswap"
$ws.Range("C15").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.43"
$ws.Range("E15").Value = "  -3.87%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.885.86"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.48"
$ws.Range("E18").Value = "  -3.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.11"
$ws.Range("E19").Value = "  -4.12%  "
$ws.Range("D20").Value = "67.482.52"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "426.77"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.25"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.12"
$ws.Range("E23").Value = "  -4.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.93"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.56"
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.42"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.92"
$ws.Range("E27").Value = "  +3.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "724.33"
$ws.Range("E29").Value = "  -3.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.06"
$ws.Range("E30").Value = "  -5.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.125"
$ws.Range("E31").Value = "  -6.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.79"
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "62.17"
$ws.Range("E33").Value = "  +8.34%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.34"
$ws.Range("E34").Value = "  -5.11%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0844"
$ws.Range("E35").Value = "  +17.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.149"
$ws.Range("E36").Value = "  -5.10%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.31"
$ws.Range("E38").Value = "  -3.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0460"
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "30.41"
$ws.Range("E40").Value = "  +20.97%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.72"
$ws.Range("E41").Value = "  +9.09%  "
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.97"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.89"
$ws.Range("E43").Value = "  +7.47%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.01"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.334"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.138"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.33"
$ws.Range("E47").Value = "  -3.99%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.12"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.14"
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.90"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.81"
$ws.Range("E51").Value = "  -3.14%  "
